# Apply Map125 scene update: add "Sina" / "Rewrite     -   Lily" translation
# helper cells in column C, and duplicate the "cake" label into column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C43").Value = "Sina"
$ws.Range("C44").Value = "Rewrite     -   Lily"
$ws.Range("D47").Value = "cake"
